$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new parts rows at the bottom of the list, continuing the
# existing "=<prev>+1" numbering pattern and shared-string text values.
$ws.Range("A53").Formula = "=A52+1"
$ws.Range("B53").Value = "REV1 steering rack assembly"

$ws.Range("A54").Formula = "=A53+1"
$ws.Range("B54").Value = "REV1 pedal assembly"

# Keep the active selection in sync with the newly-added last row, like
# Excel would after typing into B54.
$ws.Range("B54").Select()

# Minimize the workbook window (matches workbookView minimized="1").
$wb.Windows.Item(1).WindowState = -4140
